$wb = $excel.ActiveWorkbook

# Configurable zero_before_threshold parameter: recalculated First_Noticeable_Increase_Index,
# First_Noticeable_Increase_Cumulative_Value, and Pulse_Width for each Step3_DataPts sheet
# to allow dims before the noise_threshold / First Rise Point to be zeroed out.

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001048487683236386
$ws.Range("G2").Value = 36
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.02303314456298656
$ws.Range("G3").Value = 20
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.007764583867857662
$ws.Range("G4").Value = 36
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.02721509362612558
$ws.Range("G5").Value = 32
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.0252262077297413
$ws.Range("G6").Value = 36

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001048487683236386
$ws.Range("G2").Value = 56
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.02303314456298656
$ws.Range("G3").Value = 56
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.007764583867857662
$ws.Range("G4").Value = 56
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.02721509362612558
$ws.Range("G5").Value = 56
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.0252262077297413
$ws.Range("G6").Value = 56

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001048487683236386
$ws.Range("G2").Value = 66
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.02303314456298656
$ws.Range("G3").Value = 66
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.007764583867857662
$ws.Range("G4").Value = 65
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.02721509362612558
$ws.Range("G5").Value = 66
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.0252262077297413
$ws.Range("G6").Value = 63

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001048487683236386
$ws.Range("G2").Value = 78
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.02303314456298656
$ws.Range("G3").Value = 75
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.007764583867857662
$ws.Range("G4").Value = 83
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.02721509362612558
$ws.Range("G5").Value = 83
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.0252262077297413
$ws.Range("G6").Value = 83

Write-Host "Applied zero_before_threshold updates to Step3_DataPts sheets"
